# Biblioteca.xlsx - rendererPrestamos "issue picking titulo" fix
#
# USUARIOS sheet: the used range shrinks from A1:D10 down to A1:D2 (rows 3-10
# removed) and the one remaining data row (row 2) gets its placeholder "1"
# values bumped to "2".
#
# CATALOGO sheet: gains a third data row (A3:E3), filled with the same
# placeholder "1" values as row 2, growing the used range from A1:E2 to
# A1:E3.

$wb = $excel.ActiveWorkbook

# --- USUARIOS (sheet1.xml) ---------------------------------------------
$wsUsuarios = $wb.Worksheets.Item("USUARIOS")

# Drop rows 3-10 entirely so the sheet's used range becomes A1:D2.
$wsUsuarios.Rows("3:10").Delete()

# Row 2 values move from "1" to "2". Use a leading apostrophe so the
# digits are stored as text (matching the sheet's existing text-as-number
# cells) instead of being coerced into numeric cells.
$wsUsuarios.Range("A2").Value = "'2"
$wsUsuarios.Range("B2").Value = "'2"
$wsUsuarios.Range("C2").Value = "'2"
$wsUsuarios.Range("D2").Value = "'2"

# --- CATALOGO (sheet2.xml) ----------------------------------------------
$wsCatalogo = $wb.Worksheets.Item("CATALOGO")

# Append a new row 3 (A3:E3) with the same text "1" placeholder values used
# in row 2, growing the used range to A1:E3.
$wsCatalogo.Range("A3").Value = "'1"
$wsCatalogo.Range("B3").Value = "'1"
$wsCatalogo.Range("C3").Value = "'1"
$wsCatalogo.Range("D3").Value = "'1"
$wsCatalogo.Range("E3").Value = "'1"

$wb.Save()
